# Insert one new data row at row 329 (pushing existing rows 329:463 down to 330:464)
# and populate it with the new weekly price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("329:329").Insert()

$ws.Cells.Item(329, 1).Value = 4
$ws.Cells.Item(329, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(329, 3).Value = "Los Lagos"
$ws.Cells.Item(329, 4).Value = 45027
$ws.Cells.Item(329, 5).Value = 10
$ws.Cells.Item(329, 6).Value = 100114014
$ws.Cells.Item(329, 7).Value = "Betarraga"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 1000
$ws.Cells.Item(329, 11).Value = 1100
$ws.Cells.Item(329, 12).Value = 1200
$ws.Cells.Item(329, 13).Value = 1150
$ws.Cells.Item(329, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(329, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(329, 16).Value = 230
$ws.Cells.Item(329, 17).Value = 5
$ws.Cells.Item(329, 18).Value = "Hortaliza"

$ws.Cells.Item(329, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
